$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "view some educational content on a particular topic" and
# " (Have planet facts)" were two separate runs; they become one run with
# the combined text. Re-issuing the same text through Find/Replace across
# both runs causes Word to collapse them into a single run.
$d.Content.Find.Execute(
    "view some educational content on a particular topic (Have planet facts)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "view some educational content on a particular topic (Have planet facts)",
    2)

# --- Change 2 -------------------------------------------------------------
# "content" and "(" were two separate runs; they become one run "content(".
$d.Content.Find.Execute(
    "content(",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "content(",
    2)

# --- Change 3 ---------------------------------------------------------------
# Add a new run "add " immediately before the run containing "quiz" in the
# "Extensions" list (numId 2). Insert it as its own run (via InsertXML) so it
# stays a distinct <w:r> rather than merging into the existing "quiz" run.
$quizPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r", "`n") -eq "quiz") {
        $quizPara = $para
        break
    }
}

$insertionPoint = $d.Range($quizPara.Range.Start, $quizPara.Range.Start)
$insertionPoint.InsertXML(
    '<?xml version="1.0"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:r><w:t xml:space="preserve">add </w:t></w:r></w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>')
